$wb = $excel.ActiveWorkbook

# "Input Values" sheet holds the run metadata (Created at / PK data file)
$ws = $wb.Worksheets.Item("Input Values")

# Update "Created at" timestamp
$ws.Range("C2").Value = "2022-12-12 18:30:04.186656"

# Update "PK data file" identifier
$ws.Range("B12").Value = "890c10d9-87a7-4b49-a798-3642e182b282_pk_data.tsv"

# Force a full recalculation on next load (calcPr fullCalcOnLoad="1")
$wb.ForceFullCalculation = $true
